# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 45 (pushing the existing rows
# 45-62 down to 46-63), matching the new "Provincia del Elquí" /
# 2022-01-27 entry added by this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 45..62 down to 46..63, leaving a blank row 45 to populate.
$ws.Rows.Item(45).Insert()

$ws.Cells.Item(45, 1).Value  = 8
$ws.Cells.Item(45, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(45, 3).Value  = "Coquimbo"
$ws.Cells.Item(45, 4).Value  = 44588
$ws.Cells.Item(45, 5).Value  = 4
$ws.Cells.Item(45, 6).Value  = 100112030
$ws.Cells.Item(45, 7).Value  = "Poroto granado"
$ws.Cells.Item(45, 8).Value  = "Sin especificar"
$ws.Cells.Item(45, 9).Value  = "Primera"
$ws.Cells.Item(45, 10).Value = 500
$ws.Cells.Item(45, 11).Value = 29000
$ws.Cells.Item(45, 12).Value = 30000
$ws.Cells.Item(45, 13).Value = 29500
$ws.Cells.Item(45, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(45, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(45, 16).Value = 1180
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
